# Add support for multi language functionality.
#
# The translation sheet used to carry 6 language columns
# (key, en-gb, en-us, nl, de, el, es -> columns A..G). The "en-us" and
# "nl" columns are dropped, leaving key/en-gb/de/el/es (columns A..E).
# Deleting the two whole columns also removes the now-unreferenced
# "en-us"/"nl" shared strings and shifts every other language/text
# reference left automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C (en-us) and D (nl) are removed entirely; what was E/F/G
# (de/el/es) shifts left into C/D/E.
$ws.Range("C1:D1").EntireColumn.Delete()

# Mirror the author's final selection after the column delete.
$ws.Columns("C:C").Select() | Out-Null
